# Updates crypto price (D) and volume-change (E) columns to match the
# refreshed GitHub Actions scrape. Numeric-looking price strings are
# prefixed with a leading apostrophe so Excel stores them as literal
# text (preserving trailing zeros like '1.00') instead of coercing them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.470.32"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.454.62"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'555.82"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'160.73"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D9").Value = "2.453.83"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  -6.55%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  -5.51%  "
$ws.Range("D13").Value = "'4.75"
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").Value = "2.905.35"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "68.189.36"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -4.74%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "2.469.82"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'10.74"
$ws.Range("E19").Value = "  -3.93%  "
$ws.Range("D20").Value = "'340.29"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "'7.00"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").Value = "'3.75"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "'1.86"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "'66.25"
$ws.Range("E26").Value = "  -5.93%  "
$ws.Range("D27").Value = "2.581.13"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'8.05"
$ws.Range("E29").Value = "  -6.79%  "
$ws.Range("D30").Value = "0.0₃0812"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("E31").Value = "  -7.21%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'428.01"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("E35").Value = "  -5.31%  "
$ws.Range("D36").Value = "'155.75"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "'19.00"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("D40").Value = "'17.71"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").Value = "'4.36"
$ws.Range("E42").Value = "  -4.66%  "
$ws.Range("D43").Value = "'1.46"
$ws.Range("E43").Value = "  -7.62%  "
$ws.Range("D44").Value = "'1.09"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("E45").Value = "  -6.49%  "
$ws.Range("D46").Value = "'132.05"
$ws.Range("E46").Value = "  -4.18%  "
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'0.0712"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'0.477"
$ws.Range("E49").Value = "  -6.53%  "
$ws.Range("D50").Value = "'0.559"
$ws.Range("E50").Value = "  -2.16%  "
$ws.Range("D51").Value = "'0.0906"
$ws.Range("E51").Value = "  -1.60%  "
